$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'31.106.11"
$ws.Range("E2").Value = "  +3.78%  "
$ws.Range("D3").Value = "'1.915.64"
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'245.82"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").Value = "'0.4992"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("E8").Value = "  +3.36%  "
$ws.Range("D9").Value = "'0.06944"
$ws.Range("E9").Value = "  +5.04%  "
$ws.Range("D10").Value = "'1.914.87"
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("D11").Value = "'16.98"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "'0.07320"
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("D13").Value = "'89.66"
$ws.Range("E13").Value = "  +4.97%  "
$ws.Range("D14").Value = "'0.6831"
$ws.Range("E14").Value = "  +3.03%  "
$ws.Range("E15").Value = "  +5.17%  "
$ws.Range("D16").Value = "'31.060.66"
$ws.Range("D17").Value = "'0.000008094"
$ws.Range("E17").Value = "  +2.84%  "
$ws.Range("D18").Value = "'13.43"
$ws.Range("E18").Value = "  +5.27%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "'2.161.14"
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").Value = "'0.9985"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "'4.885"
$ws.Range("E22").Value = "  +2.80%  "
$ws.Range("D23").Value = "'175.39"
$ws.Range("E23").Value = "  +29.02%  "
$ws.Range("D24").Value = "'6.081"
$ws.Range("E24").Value = "  +9.44%  "
$ws.Range("D25").Value = "'9.335"
$ws.Range("E25").Value = "  +2.19%  "
$ws.Range("D26").Value = "'151.79"
$ws.Range("E26").Value = "  +2.28%  "
$ws.Range("D27").Value = "'18.11"
$ws.Range("E27").Value = "  +8.07%  "
$ws.Range("D28").Value = "'1.951"
$ws.Range("E28").Value = "  +2.14%  "
$ws.Range("E29").Value = "  +2.39%  "
$ws.Range("D30").Value = "'4.364"
$ws.Range("E30").Value = "  +4.46%  "
$ws.Range("D31").Value = "'0.08945"
$ws.Range("E31").Value = "  +3.62%  "
$ws.Range("D32").Value = "'4.064"
$ws.Range("E32").Value = "  +3.16%  "
$ws.Range("D33").Value = "'0.05254"
$ws.Range("E33").Value = "  +5.82%  "
$ws.Range("D34").Value = "'0.7487"
$ws.Range("E34").Value = "  +6.66%  "
$ws.Range("D36").Value = "'2.663"
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("D37").Value = "'0.01915"
$ws.Range("E37").Value = "  +16.94%  "
$ws.Range("D38").Value = "'2.741"
$ws.Range("E38").Value = "  +2.18%  "
$ws.Range("D39").Value = "'2.201"
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("D40").Value = "'0.9428"
$ws.Range("E40").Value = "  +1.46%  "
$ws.Range("D41").Value = "'5.955"
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("E42").Value = "  +4.51%  "
$ws.Range("D43").Value = "'105.23"
$ws.Range("E43").Value = "  +3.45%  "
$ws.Range("D44").Value = "'7.826"
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").Value = "'0.1332"
$ws.Range("E46").Value = "  +5.94%  "
$ws.Range("D47").Value = "'0.05858"
$ws.Range("E47").Value = "  +2.61%  "
$ws.Range("D48").Value = "'8.628"
$ws.Range("E48").Value = "  +4.69%  "
$ws.Range("D49").Value = "'0.3897"
$ws.Range("E49").Value = "  +5.46%  "
$ws.Range("D50").Value = "'33.36"
$ws.Range("E50").Value = "  +2.78%  "
$ws.Range("E51").Value = "  +4.51%  "
